$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Build the (large) stack-trace text that will become the new "Actual" value
# on the "Case 1" sheet. Single-quoted here-string => no variable expansion,
# so literal '$' and "'" characters inside the trace are kept verbatim.
# ---------------------------------------------------------------------------
$stackTrace = @'
java.lang.Exception: False to find 'Xác nhận/Text'!
	at helper.ElementHelper.findElementByAccessibilityId(ElementHelper.java:77)
	at reader.CashinReader.read(CashinReader.java:53)
	at reader.CashinReader.read(CashinReader.java:8)
	at com.automation.test.reader.AbstractReader.getData(AbstractReader.java:15)
	at com.automation.test.TestVerification.getActual(TestVerification.java:85)
	at com.automation.test.TestVerification.isOk(TestVerification.java:122)
	at com.automation.test.TestAction.run(TestAction.java:85)
	at com.automation.test.TestCase.run(TestCase.java:43)
	at FirstAndroidTest.enterForm(FirstAndroidTest.java:53)
	at java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:104)
	at java.base/java.lang.reflect.Method.invoke(Method.java:577)
	at org.testng.internal.invokers.MethodInvocationHelper.invokeMethod(MethodInvocationHelper.java:139)
	at org.testng.internal.invokers.TestInvoker.invokeMethod(TestInvoker.java:677)
	at org.testng.internal.invokers.TestInvoker.invokeTestMethod(TestInvoker.java:221)
	at org.testng.internal.invokers.MethodRunner.runInSequence(MethodRunner.java:50)
	at org.testng.internal.invokers.TestInvoker$MethodInvocationAgent.invoke(TestInvoker.java:962)
	at org.testng.internal.invokers.TestInvoker.invokeTestMethods(TestInvoker.java:194)
	at org.testng.internal.invokers.TestMethodWorker.invokeTestMethods(TestMethodWorker.java:148)
	at org.testng.internal.invokers.TestMethodWorker.run(TestMethodWorker.java:128)
	at java.base/java.util.ArrayList.forEach(ArrayList.java:1511)
	at org.testng.TestRunner.privateRun(TestRunner.java:806)
	at org.testng.TestRunner.run(TestRunner.java:601)
	at org.testng.SuiteRunner.runTest(SuiteRunner.java:433)
	at org.testng.SuiteRunner.runSequentially(SuiteRunner.java:427)
	at org.testng.SuiteRunner.privateRun(SuiteRunner.java:387)
	at org.testng.SuiteRunner.run(SuiteRunner.java:330)
	at org.testng.SuiteRunnerWorker.runSuite(SuiteRunnerWorker.java:52)
	at org.testng.SuiteRunnerWorker.run(SuiteRunnerWorker.java:95)
	at org.testng.TestNG.runSuitesSequentially(TestNG.java:1256)
	at org.testng.TestNG.runSuitesLocally(TestNG.java:1176)
	at org.testng.TestNG.runSuites(TestNG.java:1099)
	at org.testng.TestNG.run(TestNG.java:1067)
	at com.intellij.rt.testng.IDEARemoteTestNG.run(IDEARemoteTestNG.java:66)
	at com.intellij.rt.testng.RemoteTestNGStarter.main(RemoteTestNGStarter.java:109)
Caused by: org.openqa.selenium.WebDriverException: An unknown server-side error occurred while processing the command. Original error: 'POST /elements' cannot be proxied to UiAutomator2 server because the instrumentation process is not running (probably crashed). Check the server log and/or the logcat output for more details
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'LONGTRAN2-ITC-C', ip: '10.40.128.103', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '18.0.1.1'
Driver info: io.appium.java_client.AppiumDriver
Capabilities {app: D:/JavaCode/app/MOMO.apk, appPackage: vn.momo.platform.test, automationName: Appium, databaseEnabled: false, desired: {app: D:/JavaCode/app/MOMO.apk, automationName: Appium, deviceName: 127.0.0.1:21503, noReset: true, platformName: android, platformVersion: 7.1.2}, deviceApiLevel: 25, deviceManufacturer: google, deviceModel: G011A, deviceName: 127.0.0.1:21503, deviceScreenDensity: 300, deviceScreenSize: 900x1600, deviceUDID: 127.0.0.1:21503, javascriptEnabled: true, locationContextEnabled: false, networkConnectionEnabled: true, noReset: true, pixelRatio: 1.875, platform: LINUX, platformName: Android, platformVersion: 7.1.2, statBarHeight: 45, takesScreenshot: true, viewportRect: {height: 1555, left: 0, top: 45, width: 900}, warnings: {}, webStorageEnabled: false}
Session ID: 3a10473a-6acf-4958-8705-22c3a629f444
*** Element info: {Using=accessibility id, value=Xác nhận/Text}
	at java.base/jdk.internal.reflect.DirectConstructorHandleAccessor.newInstance(DirectConstructorHandleAccessor.java:67)
	at java.base/java.lang.reflect.Constructor.newInstanceWithCaller(Constructor.java:499)
	at java.base/java.lang.reflect.Constructor.newInstance(Constructor.java:483)
	at org.openqa.selenium.remote.http.W3CHttpResponseCodec.createException(W3CHttpResponseCodec.java:187)
	at org.openqa.selenium.remote.http.W3CHttpResponseCodec.decode(W3CHttpResponseCodec.java:122)
	at org.openqa.selenium.remote.http.W3CHttpResponseCodec.decode(W3CHttpResponseCodec.java:49)
	at org.openqa.selenium.remote.HttpCommandExecutor.execute(HttpCommandExecutor.java:158)
	at io.appium.java_client.remote.AppiumCommandExecutor.execute(AppiumCommandExecutor.java:250)
	at org.openqa.selenium.remote.RemoteWebDriver.execute(RemoteWebDriver.java:552)
	at io.appium.java_client.DefaultGenericMobileDriver.execute(DefaultGenericMobileDriver.java:41)
	at io.appium.java_client.AppiumDriver.execute(AppiumDriver.java:1)
	at org.openqa.selenium.remote.RemoteWebDriver.findElements(RemoteWebDriver.java:353)
	at io.appium.java_client.DefaultGenericMobileDriver.findElements(DefaultGenericMobileDriver.java:53)
	at io.appium.java_client.AppiumDriver.findElements(AppiumDriver.java:179)
	at io.appium.java_client.FindsByAccessibilityId.findElementsByAccessibilityId(FindsByAccessibilityId.java:50)
	at io.appium.java_client.AppiumDriver.findElementsByAccessibilityId(AppiumDriver.java:217)
	at helper.ElementHelper.findElementByAccessibilityIdIsDisplayed(ElementHelper.java:19)
	at helper.ElementHelper.findElementByAccessibilityId(ElementHelper.java:73)
	... 33 more
'@
# Re-attach the trailing newline the original shared string carries
# (here-strings drop the very last line break before '@).
$stackTrace = $stackTrace + "`n"

$colorRed   = 255      # RGB(255,0,0)   -> FAILED
$colorGreen = 32768    # RGB(0,128,0)   -> PASSED

# ---------------------------------------------------------------------------
# Summary sheet: report time, pass/fail counters and the test-case row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B1").Value = "2022-06-07T18:01:30.988113500"

$summary.Range("B3").Value = 0.0   # Passed count
$summary.Range("B4").Value = 1.0   # Failed count

$summary.Range("C9").Value = "FAILED"
$summary.Range("D9").Value = "00:00:59.768"

$summaryRow9 = $summary.Range("A9:E9")
$summaryRow9.Font.Color = $colorRed

# ---------------------------------------------------------------------------
# "Case 1" sheet: result cell & the actual-outcome (now an exception trace).
# ---------------------------------------------------------------------------
$case1 = $wb.Worksheets.Item("Case 1")

$case1.Range("E3").Value = "FAILED"
$case1.Range("E3").Font.Color = $colorRed

$case1.Range("F3").Value = $stackTrace
# Setting a long, multi-line value can make the engine auto-grow the row;
# restore the row to its natural (default) height so no stray "ht"/
# "customHeight" attribute is introduced, matching the original layout.
$case1.Rows.Item(3).EntireRow.AutoFit()
